# feat: added custom factor feature, plus simple UI
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update TYPE on row 2 from "Indicators" to "Custom"
$ws.Range("B2").Value = "Custom"

# Update EPOCHS value on row 2 from 300 to 200
$ws.Range("E2").Value = 200

# Clear the second data row (row 3), keeping formatting/styles
$ws.Range("A3:G3").ClearContents()

# Update the active selection to E2
$ws.Range("E2").Select()

$wb.Save()
